$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text so number-like values
# (e.g. "10.15", "523.84") are not auto-converted to numeric cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.621.04"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "3.846.42"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "523.84"
$ws.Range("E5").Value = "  +7.45%  "

$ws.Range("D6").Value = "143.19"
$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("D7").Value = "0.605"
$ws.Range("E7").Value = "  -2.78%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "0.710"
$ws.Range("E9").Value = "  -4.53%  "

$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  -6.15%  "

$ws.Range("D11").Value = "0.0000325"
$ws.Range("E11").Value = "  -7.49%  "

$ws.Range("D12").Value = "41.60"
$ws.Range("E12").Value = "  -3.29%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "10.15"
$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.451.26"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "21.38"
$ws.Range("E15").Value = "  +6.82%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.853.37"
$ws.Range("E16").Value = "  -1.82%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "1.23"
$ws.Range("E17").Value = "  +6.74%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "13.87"
$ws.Range("E18").Value = "  -2.50%  "

$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("D20").Value = "68.683.27"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").Value = "418.79"
$ws.Range("E21").Value = "  -3.52%  "

$ws.Range("D22").Value = "3.40"
$ws.Range("E22").Value = "  -3.56%  "

$ws.Range("D23").Value = "13.97"
$ws.Range("E23").Value = "  -5.08%  "

$ws.Range("D24").Value = "86.45"
$ws.Range("E24").Value = "  -4.79%  "

$ws.Range("D25").Value = "3.96"
$ws.Range("E25").Value = "  +5.53%  "

$ws.Range("D26").Value = "11.29"
$ws.Range("E26").Value = "  -9.44%  "

$ws.Range("D27").Value = "10.52"
$ws.Range("E27").Value = "  -4.40%  "

$ws.Range("D28").Value = "35.92"
$ws.Range("E28").Value = "  -3.89%  "

$ws.Range("D29").Value = "680.39"
$ws.Range("E29").Value = "  -4.74%  "

$ws.Range("D30").Value = "13.06"
$ws.Range("E30").Value = "  -2.37%  "

$ws.Range("D31").Value = "0.125"
$ws.Range("E31").Value = "  -3.89%  "

$ws.Range("D32").Value = "2.83"
$ws.Range("E32").Value = "  -3.07%  "

$ws.Range("D33").Value = "67.13"
$ws.Range("E33").Value = "  +9.46%  "

$ws.Range("D34").Value = "0.436"
$ws.Range("E34").Value = "  +5.26%  "

$ws.Range("D35").Value = "0.0₃0847"
$ws.Range("E35").Value = "  -3.92%  "

$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  -4.00%  "

$ws.Range("D37").Value = "39.58"
$ws.Range("E37").Value = "  -3.33%  "

$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("D39").Value = "0.146"

$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0475"
$ws.Range("E41").Value = "  -3.61%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "3.15"
$ws.Range("E42").Value = "  +1.71%  "

$ws.Range("D43").Value = "3.14"
$ws.Range("E43").Value = "  +2.37%  "

$ws.Range("D44").Value = "2.71"
$ws.Range("E44").Value = "  -8.58%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "0.139"
$ws.Range("E46").Value = "  -2.97%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.762.50"
$ws.Range("E47").Value = "  +14.54%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.93"
$ws.Range("E48").Value = "  +4.30%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000271"
$ws.Range("E49").Value = "  +11.54%  "

$ws.Range("D50").Value = "143.89"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").Value = "3.25"
$ws.Range("E51").Value = "  -3.47%  "
